# Update cryptos list - GitHub Actions scheduled data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row => @{ D = newPriceString; E = newVolumeString }
$updates = @{
    2  = @{ D = "67.432.43"; E = "  +0.89%  " }
    3  = @{ D = "3.494.35" }
    4  = @{ E = "  -0.02%  " }
    5  = @{ D = "598.44"; E = "  +0.82%  " }
    6  = @{ D = "180.16"; E = "  +4.46%  " }
    7  = @{ D = "0.608"; E = "  +4.39%  " }
    8  = @{ D = "1.00"; E = "  +0.02%  " }
    9  = @{ D = "3.495.46"; E = "  +0.04%  " }
    11 = @{ E = "  -1.81%  " }
    12 = @{ E = "  +1.42%  " }
    13 = @{ D = "4.099.29"; E = "  +0.04%  " }
    14 = @{ D = "32.31"; E = "  +10.32%  " }
    15 = @{ E = "  +0.38%  " }
    16 = @{ D = "67.412.27"; E = "  +0.86%  " }
    17 = @{ E = "  -0.02%  " }
    18 = @{ D = "3.494.85"; E = "  -0.04%  " }
    19 = @{ E = "  +0.50%  " }
    20 = @{ D = "14.29"; E = "  +0.63%  " }
    21 = @{ D = "390.24"; E = "  -0.79%  " }
    22 = @{ D = "7.95"; E = "  +0.32%  " }
    23 = @{ D = "73.96"; E = "  +1.02%  " }
    24 = @{ D = "0.542"; E = "  +1.61%  " }
    25 = @{ D = "0.999"; E = "  -0.03%  " }
    26 = @{ D = "5.74" }
    27 = @{ E = "  +0.60%  " }
    28 = @{ E = "  +1.63%  " }
    29 = @{ E = "  -2.78%  " }
    30 = @{ E = "  +0.24%  " }
    31 = @{ D = "6.19"; E = "  +0.98%  " }
    32 = @{ E = "  +0.22%  " }
    33 = @{ E = "  +1.12%  " }
    34 = @{ D = "23.55"; E = "  -0.35%  " }
    35 = @{ D = "7.40"; E = "  +0.70%  " }
    36 = @{ E = "  +0.02%  " }
    37 = @{ E = "  -0.56%  " }
    38 = @{ D = "163.33"; E = "  +0.42%  " }
    39 = @{ E = "  -0.63%  " }
    40 = @{ E = "  +11.11%  " }
    41 = @{ E = "  -0.76%  " }
    42 = @{ D = "6.85"; E = "  +0.05%  " }
    43 = @{ D = "4.64"; E = "  +0.29%  " }
    44 = @{ D = "2.852.64"; E = "  +0.38%  " }
    45 = @{ D = "26.51"; E = "  +1.75%  " }
    46 = @{ D = "26.86"; E = "  -1.23%  " }
    47 = @{ D = "0.0724"; E = "  -1.89%  " }
    48 = @{ D = "41.70" }
    49 = @{ E = "  -0.19%  " }
    50 = @{ D = "333.58"; E = "  -1.42%  " }
    51 = @{ E = "  -1.33%  " }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $rowData["D"]
    }
    if ($rowData.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $rowData["E"]
    }
}
